# solve for 18R first in calculate_17R.py
# Updates the recalculated 15Rbulk (F,L), 17R (G,M), gamma (N) and kappa (O)
# columns on both worksheets ("all" and "ATM-B6") for data rows 2-8.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = @{ F = 0.003702222127051057; G = 0.002099893460859661; L = 0.003677865467017766; M = 0.00209340517209855;  N = 0.3670017863326024; O = -0.1907152885455783 }
    3 = @{ F = 0.003702222127051057; G = 0.002099893460859661; L = 0.003673082346325243; M = 0.002089664688983711; N = 0.379263849191313;   O = -0.1782194751156673 }
    4 = @{ F = 0.003697034627764989; G = 0.002096007444871974; L = 0.003677865467017766; M = 0.00209340517209855;  N = 0.2856629739041236; O = -0.2720332900947778 }
    5 = @{ F = 0.003697034627764989; G = 0.002096007444871974; L = 0.003673082346325243; M = 0.002089664688983711; N = 0.2979250367628342; O = -0.2595374766648668 }
    6 = @{ F = 0.003699619699356788; G = 0.002095496242379899; L = 0.003677057440048185; M = 0.00209286825816626;  N = 0.2584883145367431; O = -0.2986551597096462 }
    7 = @{ F = 0.003696192533421012; G = 0.002095322040949255; L = 0.003677057440048185; M = 0.00209286825816626;  N = 0.3694649080647032; O = -0.1877069600126482 }
    8 = @{ F = 0.003699194886987533; G = 0.002097478198739332; L = 0.003677057440048185; M = 0.00209286825816626;  N = 0.3028353265431929; O = -0.254319494073837  }
}

foreach ($sheetName in @("all", "ATM-B6")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $cols = $newValues[$row]
        $ws.Range("F$row").Value = $cols.F
        $ws.Range("G$row").Value = $cols.G
        $ws.Range("L$row").Value = $cols.L
        $ws.Range("M$row").Value = $cols.M
        $ws.Range("N$row").Value = $cols.N
        $ws.Range("O$row").Value = $cols.O
    }
}
